$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.217382907867432
$ws.Range("B1").Value = 5.306960105895996
$ws.Range("C1").Value = 4.38273811340332
$ws.Range("D1").Value = 5.148453712463379
$ws.Range("E1").Value = 5.188767433166504
